$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 263
$ws.Range("I2").Value = 732
$ws.Range("J2").Value = 3208
$ws.Range("K2").Value = 23
$ws.Range("L2").Value = 814
$ws.Range("M2").Value = 56
$ws.Range("N2").Value = 529
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 33
$ws.Range("S2").Value = 323
$ws.Range("T2").Value = 568
$ws.Range("U2").Value = 37
$ws.Range("V2").Value = 4796
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4658
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 72
$ws.Range("AA2").Value = 32
